$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "0.9984") are not
# auto-converted to numbers by Excel, matching the original inline-string cell type.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '29.350.73'
$ws.Range('E2').Value = '  +0.18%  '

$ws.Range('D3').Value = '1.871.20'
$ws.Range('E3').Value = '  +0.47%  '

$ws.Range('D4').Value = '0.9984'
$ws.Range('E4').Value = '  -0.22%  '

$ws.Range('D5').Value = '0.7155'
$ws.Range('E5').Value = '  +2.10%  '

$ws.Range('D6').Value = '238.39'
$ws.Range('E6').Value = '  +0.16%  '

$ws.Range('D7').Value = '0.9983'
$ws.Range('E7').Value = '  -0.25%  '

$ws.Range('D8').Value = '0.07897'
$ws.Range('E8').Value = '  -3.90%  '

$ws.Range('D9').Value = '0.3077'
$ws.Range('E9').Value = '  +1.02%  '

$ws.Range('D10').Value = '25.46'
$ws.Range('E10').Value = '  +9.15%  '

$ws.Range('E11').Value = '  +0.06%  '

$ws.Range('D12').Value = '1.870.58'
$ws.Range('E12').Value = '  +0.63%  '

$ws.Range('D13').Value = '5.249'
$ws.Range('E13').Value = '  +1.31%  '

$ws.Range('D14').Value = '0.7238'
$ws.Range('E14').Value = '  +0.82%  '

$ws.Range('D15').Value = '89.52'
$ws.Range('E15').Value = '  +0.14%  '

$ws.Range('D16').Value = '29.415.13'
$ws.Range('E16').Value = '  +0.33%  '

$ws.Range('D17').Value = '5.844'
$ws.Range('E17').Value = '  +1.02%  '

$ws.Range('D18').Value = '242.67'
$ws.Range('E18').Value = '  +2.04%  '

$ws.Range('D19').Value = '0.000007816'
$ws.Range('E19').Value = '  -0.86%  '

$ws.Range('D20').Value = '13.30'
$ws.Range('E20').Value = '  -0.82%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.120.85'
$ws.Range('E21').Value = '  +0.28%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '0.9992'
$ws.Range('E22').Value = '  -0.08%  '

$ws.Range('D23').Value = '0.9984'
$ws.Range('E23').Value = '  -0.28%  '

$ws.Range('D24').Value = '7.621'
$ws.Range('E24').Value = '  +2.01%  '

$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D25').Value = '0.1468'
$ws.Range('E25').Value = '  +1.77%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '161.78'
$ws.Range('E26').Value = '  -0.33%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '8.963'
$ws.Range('E27').Value = '  -0.40%  '

$ws.Range('D28').Value = '18.19'
$ws.Range('E28').Value = '  +0.31%  '

$ws.Range('D29').Value = '1.922'
$ws.Range('E29').Value = '  -3.14%  '

$ws.Range('D30').Value = '1.365'
$ws.Range('E30').Value = '  -4.67%  '

$ws.Range('D31').Value = '1.481'
$ws.Range('E31').Value = '  -0.37%  '

$ws.Range('D32').Value = '4.345'
$ws.Range('E32').Value = '  -2.01%  '

$ws.Range('D33').Value = '4.065'
$ws.Range('E33').Value = '  +0.08%  '

$ws.Range('D34').Value = '0.05226'
$ws.Range('E34').Value = '  +0.14%  '

$ws.Range('D35').Value = '1.189'
$ws.Range('E35').Value = '  +1.34%  '

$ws.Range('D36').Value = '0.7212'
$ws.Range('E36').Value = '  +2.24%  '

$ws.Range('D37').Value = '1.009'
$ws.Range('E37').Value = '  +0.30%  '

$ws.Range('D38').Value = '2.672'
$ws.Range('E38').Value = '  +0.30%  '

$ws.Range('D39').Value = '0.01858'
$ws.Range('E39').Value = '  +0.23%  '

$ws.Range('D40').Value = '2.704'
$ws.Range('E40').Value = '  -0.68%  '

$ws.Range('D41').Value = '1.184.93'
$ws.Range('E41').Value = '  +3.43%  '

$ws.Range('D42').Value = '0.9185'
$ws.Range('E42').Value = '  -0.34%  '

$ws.Range('D43').Value = '6.009'
$ws.Range('E43').Value = '  +0.52%  '

$ws.Range('D44').Value = '0.4297'
$ws.Range('E44').Value = '  +0.34%  '

$ws.Range('D45').Value = '71.42'
$ws.Range('E45').Value = '  +0.79%  '

$ws.Range('D46').Value = '0.9978'
$ws.Range('E46').Value = '  -0.23%  '

$ws.Range('E47').Value = '  -0.66%  '

$ws.Range('D48').Value = '0.5356'
$ws.Range('E48').Value = '  -1.23%  '

$ws.Range('D49').Value = '1.767'
$ws.Range('E49').Value = '  -0.55%  '

$ws.Range('D50').Value = '9.240'
$ws.Range('E50').Value = '  +0.53%  '

$ws.Range('D51').Value = '7.046'
$ws.Range('E51').Value = '  +0.93%  '
